$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - add row 4 for the new handback file
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-09-06 17:19:36"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2031f0a2b1a493f893ff63f1f1858cef2031f0/e2e/f2031f0a-2b1a-493f-893f-f63f1f1858ce.md",
    $null,
    $null,
    "e2e\f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
) | Out-Null

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - add row 4 for the new handback file
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-06 17:19:31"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$wsZhCn.Range("J4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-09-06 17:19:49"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2031f0a2b1a493f893ff63f1f1858cef2031f0/e2e/f2031f0a-2b1a-493f-893f-f63f1f1858ce.md",
    $null,
    $null,
    "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f2031f0a2b1a493f893ff63f1f1858cef2031f0/e2e/f2031f0a-2b1a-493f-893f-f63f1f1858ce.md",
    $null,
    $null,
    "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
) | Out-Null

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - add row 4 for the new handback file
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-06 17:19:36"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$wsDeDe.Range("J4").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-09-06 17:19:58"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2031f0a2b1a493f893ff63f1f1858cef2031f0/e2e/f2031f0a-2b1a-493f-893f-f63f1f1858ce.md",
    $null,
    $null,
    "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f2031f0a2b1a493f893ff63f1f1858cef2031f0/e2e/f2031f0a-2b1a-493f-893f-f63f1f1858ce.md",
    $null,
    $null,
    "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
) | Out-Null

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P4"))
